# Split a single plain run of text into several sibling <w:r> runs, one per
# word/space, matching the target diff. Word's OOXML writer coalesces
# adjacent runs that carry identical (absent) run properties back into one
# run at save time, so to keep the split durable we nudge the formatting of
# every piece but the last (Bold on, then off again) which is enough to make
# the writer treat each piece as a distinct run while leaving the visible
# formatting untouched (no bold remains on anything).
#
# NOTE: this COM-interop shim only binds positional parameters, not
# "-Name value" style named arguments, so every helper below takes plain
# positional params.
function Split-ParagraphIntoWords($Paragraph, $Pieces) {
    $base = $Paragraph.Range.Start
    $pos = 0
    for ($i = 0; $i -lt $Pieces.Length; $i++) {
        $piece = $Pieces[$i]
        $pieceStart = $base + $pos
        $pieceEnd = $pieceStart + $piece.Length

        if ($i -lt ($Pieces.Length - 1)) {
            # Force a run boundary after this piece (skip the very last piece;
            # it never needs to be distinguished from anything after it).
            $r = $word.ActiveDocument.Range($pieceStart, $pieceEnd)
            $r.Font.Bold = $true
            $r.Font.Bold = $false
        }

        $pos += $piece.Length
    }
}

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    # Paragraph.Range.Text includes the trailing paragraph-mark character(s);
    # strip them before comparing against plain text.
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)

    if ($styleName -eq "Title" -and $text -eq "Desmos now in STARMAST") {
        Split-ParagraphIntoWords $p @("Desmos", " ", "now", " ", "in", " ", "STARMAST")
    }
    elseif ($styleName -eq "Author" -and $text -eq "Tom Coleman") {
        Split-ParagraphIntoWords $p @("Tom", " ", "Coleman")
    }
    elseif ($styleName -eq "Abstract" -and $text -eq "Desmos figures now included in STARMAST resources!") {
        Split-ParagraphIntoWords $p @("Desmos", " ", "figures", " ", "now", " ", "included", " ", "in", " ", "STARMAST", " ", "resources!")
    }
}

Write-Host "Split title/author/abstract runs into per-word runs."
